$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1, matching the formatting of the
# other header cells (e.g. G1 - bold, centered, bordered header style).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the numeric value for the new "Save" column in row 2.
$ws.Range("H2").Value = 0
